$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 600
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 666.6667
$ws.Range("K32").Value = 400
$ws.Range("L32").Value = 666.6667
$ws.Range("M32").Value = -74
$ws.Range("N32").Value = -1318.6667

$ws.Range("H80").Value = 603.2941
$ws.Range("I80").Value = 740
$ws.Range("J80").Value = 546.3333
$ws.Range("K80").Value = 2220
$ws.Range("L80").Value = 1638.9999
$ws.Range("M80").Value = -1222
$ws.Range("N80").Value = -3634.9999

$ws.Range("H83").Value = 603.2941
$ws.Range("I83").Value = 740
$ws.Range("J83").Value = 546.3333
$ws.Range("K83").Value = 6660
$ws.Range("L83").Value = 4916.9997
$ws.Range("M83").Value = -1668
$ws.Range("N83").Value = -14900.9997

$ws.Range("H88").Value = 695176.7
$ws.Range("I88").Value = 2778703.5
$ws.Range("J88").Value = 667.75
$ws.Range("K88").Value = 2778703.5
$ws.Range("L88").Value = 667.75
$ws.Range("M88").Value = -2778297.5
$ws.Range("N88").Value = -1479.75

$ws.Range("H91").Value = 695176.7
$ws.Range("I91").Value = 2778703.5
$ws.Range("J91").Value = 667.75
$ws.Range("K91").Value = 2778703.5
$ws.Range("L91").Value = 667.75
$ws.Range("M91").Value = -2777299.5
$ws.Range("N91").Value = -3475.75

$ws.Range("H92").Value = 432.65
$ws.Range("I92").Value = 174.15384
$ws.Range("J92").Value = 912.7143
$ws.Range("K92").Value = 174.15384
$ws.Range("L92").Value = 912.7143
$ws.Range("M92").Value = 1073.84616
$ws.Range("N92").Value = -3408.7143

$ws.Range("H112").Value = 1998.425
$ws.Range("J112").Value = 2126.6758
$ws.Range("L112").Value = 6380.0274
$ws.Range("N112").Value = -8596.027399999999

$ws.Range("H132").Value = 5053054.5
$ws.Range("I132").Value = 6495357.5
$ws.Range("J132").Value = 4992.5
$ws.Range("K132").Value = 19486072.5
$ws.Range("L132").Value = 14977.5
$ws.Range("M132").Value = -19483542.5
$ws.Range("N132").Value = -20037.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 41667900
$ws.Range("J74").Value = 1878.5
$ws.Range("L74").Value = 1878.5
$ws.Range("N74").Value = -3626.5

$ws.Range("H77").Value = 41667900
$ws.Range("J77").Value = 1878.5
$ws.Range("L77").Value = 9392.5
$ws.Range("N77").Value = -18128.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 975.8889
$ws.Range("I107").Value = 975.8889
$ws.Range("K107").Value = 975.8889
$ws.Range("M107").Value = 944.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2364.95
$ws.Range("I132").Value = 1584.8667
$ws.Range("J132").Value = 4705.2
$ws.Range("K132").Value = 4754.6001
$ws.Range("L132").Value = 14115.6
$ws.Range("M132").Value = -2224.6001
$ws.Range("N132").Value = -19175.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2900.6
$ws.Range("I46").Value = 2251.5
$ws.Range("K46").Value = 6754.5
$ws.Range("M46").Value = -6663.5

$ws.Range("H57").Value = 5555.8823
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 5555.8823
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 16667.6469
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -17785.6469

$ws.Range("H58").Value = 4067.6
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 4067.6
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 12202.8
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -12458.8

$ws.Range("H59").Value = 3350
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 3350
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 10050
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -11130

$ws.Range("H60").Value = 1063.75
$ws.Range("I60").Value = 197.5
$ws.Range("J60").Value = 1930
$ws.Range("K60").Value = 592.5
$ws.Range("L60").Value = 5790
$ws.Range("M60").Value = -341.5
$ws.Range("N60").Value = -6292

$ws.Range("H61").Value = 597.35297
$ws.Range("J61").Value = 597.35297
$ws.Range("L61").Value = 1792.05891
$ws.Range("N61").Value = -2222.05891

$ws.Range("H62").Value = 2899.5
$ws.Range("I62").Value = 999
$ws.Range("J62").Value = 4800
$ws.Range("K62").Value = 2997
$ws.Range("L62").Value = 14400
$ws.Range("M62").Value = -2311
$ws.Range("N62").Value = -15772

$ws.Range("H63").Value = 4329.3335
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 4994
$ws.Range("K63").Value = 9000
$ws.Range("L63").Value = 14982
$ws.Range("M63").Value = -8251
$ws.Range("N63").Value = -16480

$ws.Range("H64").Value = 6400
$ws.Range("J64").Value = 6400
$ws.Range("L64").Value = 19200
$ws.Range("N64").Value = -19740

$ws.Range("H65").Value = 2899.5
$ws.Range("I65").Value = 999
$ws.Range("J65").Value = 4800
$ws.Range("K65").Value = 8991
$ws.Range("L65").Value = 43200
$ws.Range("M65").Value = -5559
$ws.Range("N65").Value = -50064

$ws.Range("H66").Value = 4329.3335
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 4994
$ws.Range("K66").Value = 27000
$ws.Range("L66").Value = 44946
$ws.Range("M66").Value = -23256
$ws.Range("N66").Value = -52434

$ws.Range("H67").Value = 6400
$ws.Range("J67").Value = 6400
$ws.Range("L67").Value = 19200
$ws.Range("N67").Value = -21072

$ws.Range("H131").Value = 46288.297
$ws.Range("I131").Value = 418.57144
$ws.Range("J131").Value = 54966.35
$ws.Range("K131").Value = 1255.71432
$ws.Range("L131").Value = 164899.05
$ws.Range("M131").Value = 3784.28568
$ws.Range("N131").Value = -174979.05

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 335918
$ws.Range("I14").Value = 335918
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 335918
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -335750
$ws.Range("N14").ClearContents()

$ws.Range("H122").Value = 2138.5264
$ws.Range("I122").Value = 2090.6667
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6272.000100000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3822.000100000001
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1700
$ws.Range("I100").Value = 1550
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1550
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1009
$ws.Range("N100").Value = -3082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1553.5938
$ws.Range("I132").Value = 1002.5294
$ws.Range("K132").Value = 3007.5882
$ws.Range("M132").Value = -477.5882000000001

$ws.Range("H136").Value = 3872.1025
$ws.Range("J136").Value = 9674.615
$ws.Range("L136").Value = 29023.845
$ws.Range("N136").Value = -34123.845
